# "added 7th week hw"
#
# Duplicate the "Yenidosya" worksheet (ISIM/SOYISIM/YAS/ALDIGI MAAS/CINSIYETI
# table) into a brand-new worksheet named "Yenidosya5", appended as the last
# sheet of the workbook.

$wb = $excel.ActiveWorkbook

$src = $wb.Worksheets.Item("Yenidosya")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Yenidosya5"

# Every cell in the workbook was carrying a redundant explicit "General"
# number-format style (the only style besides the default). Clear the
# formatting on all sheets (old and new) so the cells fall back to the
# workbook's default/unstyled cell format.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.ClearFormats()
}
